# This script updates the "Assignment" (TA) column (column D) of Sheet1 in the
# Capstone Project TA-matching workbook. The commit message says the backend
# changed the order of sorting, which resulted in a reshuffled TA assignment
# roster; the only semantic change is the text in column D for the affected
# rows (the Course/Section/Instructor/Course Title/Days/Times/Building/Room
# columns are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = "TA: Anurag Nandyala"
$ws.Range("D14").Value = "TA: Anurag Nandyala"
$ws.Range("D15").Value = "TA: Anurag Nandyala"
$ws.Range("D16").Value = "TA: Anurag Nandyala"
$ws.Range("D29").Value = "TA: Amir Faiyaz (DGSA), TA: Ban Tran (DGSA)"
$ws.Range("D30").Value = "TA: Amir Faiyaz (DGSA), TA: Ban Tran (DGSA), TA: Hasan Alqudah"
$ws.Range("D31").Value = "TA: Linfeng Wei, TA: Amir Faiyaz (DGSA), TA: Ban Tran (DGSA)"
$ws.Range("D32").Value = "TA: Amir Faiyaz (DGSA), TA: Ban Tran (DGSA), TA: Hasan Alqudah"
$ws.Range("D33").Value = "TA: Sthephany Rojas, TA: Hasan Alqudah"
$ws.Range("D34").Value = "TA: Hasan Alqudah"
$ws.Range("D35").Value = "TA: Imran Pinjari (DGSA)"
$ws.Range("D36").Value = "TA: Imran Pinjari (DGSA)"
$ws.Range("D37").Value = "TA: Rakib Hossain Rifat (DGSA), TA: Imran Pinjari (DGSA), TA: Juan Marcelo Gutierrez Carballo"
$ws.Range("D38").Value = "TA: Samin Dehbashi Sani, TA: Imran Pinjari (DGSA), TA: Juan Marcelo Gutierrez Carballo"
$ws.Range("D39").Value = "TA: Ruoyu Xu, TA: Xin Zhang (DGSA), TA: Juan Marcelo Gutierrez Carballo"
$ws.Range("D40").Value = "TA: Juan Marcelo Gutierrez Carballo"
$ws.Range("D41").Value = "TA: Sthephany Rojas, TA: Linfeng Wei"
$ws.Range("D44").Value = "TA: Linfeng Wei, TA: Linpeng Sun"
$ws.Range("D46").Value = "TA: Linpeng Sun, TA: Liyuan Gao (DGSA)"
$ws.Range("D47").Value = "TA: Rukayat Olapojoye, TA: Zhenyu Xu"
$ws.Range("D48").Value = "TA: Linpeng Sun"
$ws.Range("D49").Value = "TA: Muhammad Aziz Ullah, TA: Denish Otieno"
$ws.Range("D51").Value = "TA: Phornsawan Roemsri, TA: Linpeng Sun"
$ws.Range("D52").Value = "TA: Liyuan Gao (DGSA), TA: Md Mahabub Uz Zaman"
$ws.Range("D53").Value = "TA: Liyuan Gao (DGSA), TA: Md Mahabub Uz Zaman"
$ws.Range("D55").Value = "TA: Muhammad Aziz Ullah, TA: Liyuan Gao (DGSA)"
$ws.Range("D56").Value = "TA: Gaoxiang Li, TA: Nabonita Mitra "
$ws.Range("D57").Value = "TA: Gaoxiang Li, TA: Md Mahabub Uz Zaman"
$ws.Range("D58").Value = "TA: Seyed Soroush Tabadkani Avval"
$ws.Range("D59").Value = "TA: Gaoxiang Li, TA: Ruoyu Xu"
